$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 13373.875
$ws.Range("I33").Value = 19215.363
$ws.Range("K33").Value = 19215.363
$ws.Range("M33").Value = -18986.363
$ws.Range("H43").Value = 3515.3076
$ws.Range("J43").Value = 2600
$ws.Range("L43").Value = 2600
$ws.Range("N43").Value = -2738
$ws.Range("H76").Value = 12099.0625
$ws.Range("I76").Value = 24137
$ws.Range("J76").Value = 6627.273
$ws.Range("K76").Value = 24137
$ws.Range("L76").Value = 6627.273
$ws.Range("M76").Value = -23822
$ws.Range("N76").Value = -7257.273
$ws.Range("H79").Value = 12099.0625
$ws.Range("I79").Value = 24137
$ws.Range("J79").Value = 6627.273
$ws.Range("K79").Value = 24137
$ws.Range("L79").Value = 6627.273
$ws.Range("M79").Value = -23045
$ws.Range("N79").Value = -8811.273000000001
$ws.Range("H98").Value = 1184.24
$ws.Range("I98").Value = 1209.409
$ws.Range("K98").Value = 1209.409
$ws.Range("M98").Value = 288.5909999999999
$ws.Range("H122").Value = 1184.24
$ws.Range("I122").Value = 1209.409
$ws.Range("K122").Value = 3628.227
$ws.Range("M122").Value = -1178.227
$ws.Range("H132").Value = 4697.1177
$ws.Range("I132").Value = 4807.778
$ws.Range("J132").Value = 4572.625
$ws.Range("K132").Value = 14423.334
$ws.Range("L132").Value = 13717.875
$ws.Range("M132").Value = -11893.334
$ws.Range("N132").Value = -18777.875
$ws.Range("H133").Value = 78000
$ws.Range("J133").Value = 78000
$ws.Range("L133").Value = 78000
$ws.Range("N133").Value = -88120
$ws.Range("H135").Value = 904.64
$ws.Range("I135").Value = 603.9474
$ws.Range("J135").Value = 1856.8334
$ws.Range("K135").Value = 5435.5266
$ws.Range("L135").Value = 16711.5006
$ws.Range("M135").Value = -2900.5266
$ws.Range("N135").Value = -21781.5006
$ws.Range("H136").Value = 63296.668
$ws.Range("J136").Value = 63296.668
$ws.Range("L136").Value = 63296.668
$ws.Range("N136").Value = -73496.66800000001
$ws.Range("H138").Value = 6670642
$ws.Range("J138").Value = 11911089
$ws.Range("L138").Value = 35733267
$ws.Range("N138").Value = -35743547
$ws.Range("H140").Value = 199999
$ws.Range("J140").Value = 199999
$ws.Range("L140").Value = 199999
$ws.Range("N140").Value = -210359
$ws.Range("H141").Value = 2299.7693

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 20909.21
$ws.Range("I110").Value = 28522
$ws.Range("K110").Value = 28522
$ws.Range("M110").Value = -26477

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2255.5833
$ws.Range("I20").Value = 2383.625
$ws.Range("J20").Value = 1999.5
$ws.Range("K20").Value = 2383.625
$ws.Range("L20").Value = 1999.5
$ws.Range("M20").Value = -2136.625
$ws.Range("N20").Value = -2493.5
$ws.Range("H134").Value = 2904.6191
$ws.Range("I134").Value = 3077.3333
$ws.Range("K134").Value = 9231.999899999999
$ws.Range("M134").Value = -6696.999899999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H114").Value = 44215.25
$ws.Range("I114").Value = 33953.668
$ws.Range("J114").Value = 75000
$ws.Range("K114").Value = 33953.668
$ws.Range("L114").Value = 75000
$ws.Range("M114").Value = -29614.668
$ws.Range("N114").Value = -83678
$ws.Range("H132").Value = 60041.484
$ws.Range("I132").Value = 82139.88
$ws.Range("J132").Value = 4795.5
$ws.Range("K132").Value = 246419.64
$ws.Range("L132").Value = 14386.5
$ws.Range("M132").Value = -243889.64
$ws.Range("N132").Value = -19446.5
$ws.Range("H135").Value = 81195
$ws.Range("J135").Value = 81195
$ws.Range("L135").Value = 81195
$ws.Range("N135").Value = -91335
$ws.Range("H138").Value = 116615
$ws.Range("I138").Value = 38700
$ws.Range("J138").Value = 142586.67
$ws.Range("K138").Value = 38700
$ws.Range("L138").Value = 142586.67
$ws.Range("M138").Value = -33560
$ws.Range("N138").Value = -152866.67
$ws.Range("H139").Value = 126608.89
$ws.Range("I139").Value = 123283.336
$ws.Range("J139").Value = 133260
$ws.Range("K139").Value = 123283.336
$ws.Range("L139").Value = 133260
$ws.Range("M139").Value = -118143.336
$ws.Range("N139").Value = -143540
$ws.Range("H141").Value = 254443.67
$ws.Range("J141").Value = 254443.67
$ws.Range("L141").Value = 254443.67
$ws.Range("N141").Value = -264803.67

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2478.8235
$ws.Range("J68").Value = 2468.1428
$ws.Range("L68").Value = 7404.428400000001
$ws.Range("N68").Value = -9026.428400000001
$ws.Range("H71").Value = 2478.8235
$ws.Range("J71").Value = 2468.1428
$ws.Range("L71").Value = 22213.2852
$ws.Range("N71").Value = -30325.2852
$ws.Range("H107").Value = 800.4
$ws.Range("I107").Value = 799
$ws.Range("J107").Value = 800.75
$ws.Range("K107").Value = 2397
$ws.Range("L107").Value = 2402.25
$ws.Range("M107").Value = -477
$ws.Range("N107").Value = -6242.25
$ws.Range("H128").Value = 115584
$ws.Range("I128").Value = 115584
$ws.Range("K128").Value = 346752
$ws.Range("M128").Value = -341772
$ws.Range("H131").Value = 1741.9286
$ws.Range("I131").Value = 2050
$ws.Range("J131").Value = 1730.5186
$ws.Range("K131").Value = 6150
$ws.Range("L131").Value = 5191.5558
$ws.Range("M131").Value = -1110
$ws.Range("N131").Value = -15271.5558
$ws.Range("H133").Value = 6833.8335
$ws.Range("J133").Value = 20000
$ws.Range("L133").Value = 60000
$ws.Range("N133").Value = -70120

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 2523458.2
$ws.Range("I14").Value = 2523458.2
$ws.Range("K14").Value = 2523458.2
$ws.Range("M14").Value = -2523290.2
$ws.Range("H70").Value = 670002.7
$ws.Range("I70").Value = 1002504
$ws.Range("K70").Value = 1002504
$ws.Range("M70").Value = -1002234
$ws.Range("H73").Value = 670002.7
$ws.Range("I73").Value = 1002504
$ws.Range("K73").Value = 1002504
$ws.Range("M73").Value = -1001568
$ws.Range("H97").Value = 2540.524
$ws.Range("I97").Value = 1808.75
$ws.Range("K97").Value = 1808.75
$ws.Range("M97").Value = -1312.75
$ws.Range("H102").Value = 2425.5173
$ws.Range("I102").Value = 1380.2778
$ws.Range("J102").Value = 4135.909
$ws.Range("K102").Value = 1380.2778
$ws.Range("L102").Value = 4135.909
$ws.Range("M102").Value = 241.7221999999999
$ws.Range("N102").Value = -7379.909

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 2002.5
$ws.Range("I30").Value = 1931.4286
$ws.Range("J30").Value = 2500
$ws.Range("K30").Value = 1931.4286
$ws.Range("L30").Value = 2500
$ws.Range("M30").Value = -1823.4286
$ws.Range("N30").Value = -2716
$ws.Range("H46").Value = 1289.742
$ws.Range("I46").Value = 494.22726
$ws.Range("K46").Value = 494.22726
$ws.Range("M46").Value = -306.22726
$ws.Range("H82").Value = 3576.8
$ws.Range("I82").Value = 2490
$ws.Range("J82").Value = 4120.2
$ws.Range("K82").Value = 2490
$ws.Range("L82").Value = 4120.2
$ws.Range("M82").Value = -2129
$ws.Range("N82").Value = -4842.2
$ws.Range("H85").Value = 3576.8
$ws.Range("I85").Value = 2490
$ws.Range("J85").Value = 4120.2
$ws.Range("K85").Value = 2490
$ws.Range("L85").Value = 4120.2
$ws.Range("M85").Value = -1242
$ws.Range("N85").Value = -6616.2
$ws.Range("H122").Value = 4746.5713
$ws.Range("I122").Value = 3409.6
$ws.Range("J122").Value = 5281.36
$ws.Range("K122").Value = 10228.8
$ws.Range("L122").Value = 15844.08
$ws.Range("M122").Value = -7778.799999999999
$ws.Range("N122").Value = -20744.08
$ws.Range("H130").Value = 158426
$ws.Range("J130").Value = 158426
$ws.Range("L130").Value = 158426
$ws.Range("N130").Value = -168466
$ws.Range("H136").Value = 1963.8889
$ws.Range("I136").Value = 682.4400000000001
$ws.Range("K136").Value = 2047.32
$ws.Range("M136").Value = 502.6799999999998

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3918.55
$ws.Range("I126").Value = 4274.5557
$ws.Range("K126").Value = 12823.6671
$ws.Range("M126").Value = -10353.6671
$ws.Range("H136").Value = 6230.8335
$ws.Range("I136").Value = 1827
$ws.Range("J136").Value = 5481
$ws.Range("K136").Value = 5481
$ws.Range("M136").Value = -2931
